$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet: insert a new row for the "CaseHasProtectedId" setting
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")

# Insert a brand-new row 16 (pushes the old 16-32 block down by one row)
$wsSettings.Rows("16:16").Insert()

$wsSettings.Range("A16").Value = "CaseHasProtectedId"
$wsSettings.Range("B16").Value = "Skyddat Personnummer`nSkyddat personnummer"
$wsSettings.Range("B16").WrapText = $true

# Grow the sheet by one trailing (still empty) row, matching row 999's format
$wsSettings.Rows("1000:1000").RowHeight = $wsSettings.Rows("999:999").RowHeight

# ---------------------------------------------------------------------------
# Assets sheet: drop the two now-obsolete citizen-id asset rows and add the
# new "CaseHasProtectedId" asset row in their place
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Remove "SignProtectedCitizenID" (row 16) and "SignTFCitizenID" (row 17)
$wsAssets.Rows("16:17").Delete()

# The rows below shifted up by two; the row that used to be the trailing
# blank spacer (old row 20) is now row 18 - populate it with the new asset
$wsAssets.Range("A18").Value = "CaseHasProtectedId"
$wsAssets.Range("B18").Value = "CaseHasProtectedId"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved workbook state
# ---------------------------------------------------------------------------
$wsSettings.Range("B16").Select()
$wsAssets.Rows("16:16").Select()
